# Generate Report for Handoff
# Update status and timestamps across the Overview, zh-cn, and de-de sheets
# to reflect that the handoff report has been (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: B2 (zh-cn status), C2 (de-de status), D2 (Latest Handoff Date)
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-27-20 02:27:38"

# zh-cn sheet: C2 (Status), E2 (Latest Handoff Datetime)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-20 02:27:35"

# de-de sheet: C2 (Status), E2 (Latest Handoff Datetime)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-20 02:27:38"
